$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44461
$ws.Range("D55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108004
$ws.Cells.Item($row, 10).Value = "Papaya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 70
$ws.Cells.Item($row, 14).Value = 20000
$ws.Cells.Item($row, 15).Value = 21000
$ws.Cells.Item($row, 16).Value = 20429
$ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2043
$ws.Cells.Item($row, 20).Value = 10
